# Shifted interferograms to be centered about 0
#
# This edit adds a new normalisation-method comment row and a new
# (partially filled) row to the "1dmockanderrors" metadata table, and
# adds a "Comments" column to that table.
#
# Structure before the edit:
#   Row 10  : merged section header "1dmockanderrors.csv (new format)"
#   Row 11  : table header row (name, array length ... averages)
#   Row 12  : data row - 1dmockanderrors1.csv
#   Row 13  : data row - 1dmockanderrors2.csv
#
# Structure after the edit:
#   Row 11  : merged section header (shifted down by one row)
#   Row 12  : table header row (now with an extra "Comments" column)
#   Row 13  : data row - 1dmockanderrors1.csv
#   Row 14  : data row - 1dmockanderrors2.csv
#   Row 15  : NEW data row - 1dmockanderrors3.csv (with a Comments entry)
#   Row 16  : NEW (partial) row - 1dmockanderrors4.csv

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 10. This pushes the merged header,
# the table header and the two existing data rows down by one row each,
# and the second table (Table4) automatically grows/shifts with it.
$ws.Rows.Item(10).Insert()

# Grow the second table ("Table4") so it covers the new "Comments"
# column (L) and the two new data rows (15 and 16).
$lo = $ws.ListObjects.Item(2)
$lo.Resize($ws.Range("B12:L16"))

# The old J12 cell (now J13) had a custom/applied-number-format style;
# the new layout uses the plain default style for this cell.
$ws.Range("J13").Style = "Normal"
$ws.Range("J13").Value = 20

# New data row 15: 1dmockanderrors3.csv, including a comment describing
# the new normalisation method (set before the "Comments" header so the
# shared-string table is built up in the same order as the target file).
$ws.Range("L15").Value = "New normalisation meathod. Ideal (coherent, noiseless) Interferogram peaks are now set to 1."

# Header for the new "Comments" column.
$ws.Range("L12").Value = "Comments"

$ws.Range("B15").Value = "1dmockanderrors3.csv"
$ws.Range("C15").Value = 400
$ws.Range("D15").Value = 50
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 60
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 20

# New (partially filled in) row 16: just the file name so far.
$ws.Range("B16").Value = "1dmockanderrors4.csv"

# Give column L (Comments) a sensible width, matching the other
# metadata columns.
$ws.Columns.Item(12).ColumnWidth = 12

# Match the final selected cell shown in the saved workbook.
$ws.Range("C16").Select() | Out-Null
